# Generate Report for Handback
#
# Renames the first handed-back file's GUID (2568fbfd-...  -> 12c228b0-...)
# throughout the workbook, updates its handoff/handback hashes & timestamps,
# and appends a brand-new second handed-back file
# (41161395-4867-4cf8-ac67-21dd0b22646e.md) as a new row on every sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# File identifiers
# ---------------------------------------------------------------------------
$guid1 = "12c228b0-e461-436f-89e9-e53ee3f7349c"   # renamed (was 2568fbfd-...)
$guid2 = "41161395-4867-4cf8-ac67-21dd0b22646e"   # brand new second file

$hash1 = "fec4df54a5f000ea70689957276e01739547613a"   # guid1's handoff hash
$hash2 = "a8aac88ad6c0d9fd6aea61c6a0f614b5f6c511f5"   # guid2's handoff hash

# Drop every existing hyperlink on every sheet up front - the COM shim can't
# update a hyperlink's target/display text in place, so everything gets
# re-added below (both the ones that merely changed GUID and the new ones).
$ws1.Range("A1").Hyperlinks.Delete()
$ws2.Range("A1").Hyperlinks.Delete()
$ws3.Range("A1").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# Overview sheet - update row 2 (renamed file) and add row 3 (new file)
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "$guid1.md"
$ws1.Range("B2").Value = "e2e\$guid1.md"
$ws1.Range("G2").Value = "2016-08-29 17:09:09"

$ws1.Range("A3").Value = "$guid2.md"
$ws1.Range("C3").Value = ".md"
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"
$ws1.Range("G3").Value = "2016-08-29 17:09:09"
$ws1.Range("G3").NumberFormat = $dateFmt

$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c77a49f5ded47c71a809223dfbc20637c45b3884/e2e/$guid1.md", "", "", "e2e\$guid1.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c77a49f5ded47c71a809223dfbc20637c45b3884/e2e/$guid2.md", "", "", "e2e\$guid2.md") | Out-Null

$tbl1 = $ws1.ListObjects.Item("Overview")
$tbl1.Resize($ws1.Range("A1:G3"))

# ---------------------------------------------------------------------------
# zh-cn sheet - update row 2 (renamed file) and add row 3 (new file)
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "$guid1.md"
$ws2.Range("G2").Value = "$guid1.$hash1.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-29 17:08:58"
$ws2.Range("H2").NumberFormat = $dateFmt
$ws2.Range("I2").Value = "$guid1.md"
$ws2.Range("J2").Value = "$guid1.$hash1.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-29 17:09:33"
$ws2.Range("K2").NumberFormat = $dateFmt

$ws2.Range("A3").Value = "$guid2.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "True"
$ws2.Range("G3").Value = "$guid2.$hash2.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-29 17:08:58"
$ws2.Range("H3").NumberFormat = $dateFmt
$ws2.Range("I3").Value = "$guid2.md"
$ws2.Range("J3").Value = "$guid2.$hash2.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-29 17:09:33"
$ws2.Range("K3").NumberFormat = $dateFmt
$ws2.Range("M3").Value = "True"
$ws2.Range("O3").Value = "False"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c77a49f5ded47c71a809223dfbc20637c45b3884/e2e/$guid1.md", "", "", "$guid1.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3131b8a1c783b9fd4153d7d3a81f1222433f75ac/e2e/$guid1.md", "", "", "$guid1.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c77a49f5ded47c71a809223dfbc20637c45b3884/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3131b8a1c783b9fd4153d7d3a81f1222433f75ac/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null

$tbl2 = $ws2.ListObjects.Item("zh-cn")
$tbl2.Resize($ws2.Range("A1:P3"))

# ---------------------------------------------------------------------------
# de-de sheet - update row 2 (renamed file) and add row 3 (new file)
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = "$guid1.md"
$ws3.Range("G2").Value = "$guid1.$hash1.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-29 17:09:09"
$ws3.Range("H2").NumberFormat = $dateFmt
$ws3.Range("I2").Value = "$guid1.md"
$ws3.Range("J2").Value = "$guid1.$hash1.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-29 17:09:40"
$ws3.Range("K2").NumberFormat = $dateFmt

$ws3.Range("A3").Value = "$guid2.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "True"
$ws3.Range("G3").Value = "$guid2.$hash2.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-29 17:09:09"
$ws3.Range("H3").NumberFormat = $dateFmt
$ws3.Range("I3").Value = "$guid2.md"
$ws3.Range("J3").Value = "$guid2.$hash2.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-29 17:09:40"
$ws3.Range("K3").NumberFormat = $dateFmt
$ws3.Range("M3").Value = "True"
$ws3.Range("O3").Value = "False"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c77a49f5ded47c71a809223dfbc20637c45b3884/e2e/$guid1.md", "", "", "$guid1.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b8f77a09ec9892d8f5daaa08ba296681543b9321/e2e/$guid1.md", "", "", "$guid1.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c77a49f5ded47c71a809223dfbc20637c45b3884/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b8f77a09ec9892d8f5daaa08ba296681543b9321/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null

$tbl3 = $ws3.ListObjects.Item("de-de")
$tbl3.Resize($ws3.Range("A1:P3"))
